$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (title) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

# Borrow the formatted content (which carries the document's usual leading
# empty run) from an existing body paragraph so the new paragraph's
# structure matches its siblings, then swap in the real text via Find &
# Replace so the leading run survives the edit.
$donor = $d.Paragraphs(4).Range.FormattedText
$metaPara.Range.FormattedText = $donor

$metaRange = $d.Paragraphs(2).Range
$metaRange.Find.Execute(
    "Get ready to swing into the jungle with Banana King, a slot game bursting with color and charm. As you play, you'll be transported to a lush forest, where friendly chimpanzees reign supreme, munching on bananas and swinging through the trees.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description: Join the friendly chimpanzees in the lush forest and play Banana King for free with selectable pay lines, up to 60 free spins, and a special bonus feature.",
    2)

# Bold just the "Meta description" label (leave the rest, including the
# colon, un-bolded).
$metaPara = $d.Paragraphs(2)
$labelStart = $metaPara.Range.Start
$labelRange = $d.Range($labelStart, $labelStart + 16)
$labelRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# 2. Remove the old bold "Play Banana King for Free..." paragraph near the
#    end of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$oldTitlePara = $d.Paragraphs($count - 1)
$oldTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Turn the paragraph that follows it (the italic call-to-action line)
#    into the new AI image prompt, keeping the italic formatting but
#    without going through Find & Replace (which "smart quotes" the
#    straight quotes/apostrophes in the new text).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$start = $promptPara.Range.Start
$end = $promptPara.Range.End

$newText = "Prompt: Create a feature image for " + [char]34 + "Banana King" + [char]34 + " - a cartoon style image featuring a happy Maya warrior with glasses. The image should be colorful and vibrant, showcasing the lush forest full of bananas and chimpanzees. The Maya warrior should be holding a banana and smiling, with his glasses reflecting the reels of the game behind him. The background should include the intertwined branches and foliage from the game" + [char]39 + "s grid, with some of the symbols from the game as playful accents in the surrounding space. The image should convey a sense of joy and adventure, inviting players to join in on the fun of Banana King."

$insertRange = $d.Range($start, $start)
$insertRange.InsertAfter($newText)
$insertedEnd = $insertRange.End
$insertedLen = $insertedEnd - $start

# Delete the old text, which now sits right after the newly inserted text.
$oldTextStart = $insertedEnd
$oldTextEnd = $end + $insertedLen
$oldTextRange = $d.Range($oldTextStart, $oldTextEnd)
$oldTextRange.Delete()

# Re-apply italics to just the visible text (not the paragraph mark, so we
# don't introduce an extra paragraph-mark run-properties element).
$promptPara = $d.Paragraphs($d.Paragraphs.Count)
$textLen = $promptPara.Range.Text.Length
$textOnlyRange = $d.Range($promptPara.Range.Start, $promptPara.Range.Start + $textLen)
$textOnlyRange.Font.Italic = $true
